$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 (previous last row) reverts to the regular date-time format used by all other data rows
$ws.Range("A31").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 32: today's data
$ws.Cells.Item(32, 1).Value = 45616
$ws.Cells.Item(32, 2).Value = 81
$ws.Cells.Item(32, 3).Value = 66
$ws.Cells.Item(32, 4).Value = 78

# New last row (A32) takes on the distinct "last row" date-only format
$ws.Range("A32").NumberFormat = "YYYY-MM-DD"
